{"js": "// Replace the 25 \"three-digit x one-digit\" answer cells with their new\n// values. Each old value is unique in the document, and we process the\n// pairs in the same order as they appear in the document so that a value\n// which is both an \"old\" search target (index 7) and a newly produced\n// \"new\" value (index 22) never collides with itself.\nconst replacements = [\n  [\"780\u00d79=7020\", \"805\u00d77=5635\"],\n  [\"725\u00d78=5800\", \"799\u00d75=3995\"],\n  [\"246\u00d76=1476\", \"290\u00d72=580\"],\n  [\"631\u00d77=4417\", \"625\u00d75=3125\"],\n  [\"554\u00d79=4986\", \"926\u00d77=6482\"],\n  [\"860\u00d72=1720\", \"240\u00d78=1920\"],\n  [\"136\u00d76=816\", \"439\u00d72=878\"],\n  [\"917\u00d79=8253\", \"841\u00d72=1682\"],\n  [\"404\u00d77=2828\", \"464\u00d76=2784\"],\n  [\"942\u00d78=7536\", \"719\u00d78=5752\"],\n  [\"886\u00d78=7088\", \"708\u00d79=6372\"],\n  [\"243\u00d72=486\", \"535\u00d74=2140\"],\n  [\"313\u00d77=2191\", \"914\u00d75=4570\"],\n  [\"167\u00d78=1336\", \"938\u00d79=8442\"],\n  [\"797\u00d73=2391\", \"276\u00d74=1104\"],\n  [\"362\u00d79=3258\", \"933\u00d77=6531\"],\n  [\"431\u00d73=1293\", \"270\u00d73=810\"],\n  [\"447\u00d75=2235\", \"533\u00d74=2132\"],\n  [\"499\u00d76=2994\", \"181\u00d79=1629\"],\n  [\"119\u00d73=357\", \"629\u00d76=3774\"],\n  [\"970\u00d79=8730\", \"122\u00d76=732\"],\n  [\"478\u00d74=1912\", \"718\u00d79=6462\"],\n  [\"679\u00d78=5432\", \"917\u00d79=8253\"],\n  [\"804\u00d72=1608\", \"505\u00d79=4545\"],\n  [\"463\u00d73=1389\", \"431\u00d76=2586\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit x one-digit\" answer cells with their new\n# values. Each old value is unique in the document, and the pairs are\n# processed in the same order as they appear in the document so that a\n# value which is both an \"old\" search target (#8) and a newly produced\n# \"new\" value (#23) never collides with itself.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"780\u00d79=7020\", \"805\u00d77=5635\"),\n    @(\"725\u00d78=5800\", \"799\u00d75=3995\"),\n    @(\"246\u00d76=1476\", \"290\u00d72=580\"),\n    @(\"631\u00d77=4417\", \"625\u00d75=3125\"),\n    @(\"554\u00d79=4986\", \"926\u00d77=6482\"),\n    @(\"860\u00d72=1720\", \"240\u00d78=1920\"),\n    @(\"136\u00d76=816\", \"439\u00d72=878\"),\n    @(\"917\u00d79=8253\", \"841\u00d72=1682\"),\n    @(\"404\u00d77=2828\", \"464\u00d76=2784\"),\n    @(\"942\u00d78=7536\", \"719\u00d78=5752\"),\n    @(\"886\u00d78=7088\", \"708\u00d79=6372\"),\n    @(\"243\u00d72=486\", \"535\u00d74=2140\"),\n    @(\"313\u00d77=2191\", \"914\u00d75=4570\"),\n    @(\"167\u00d78=1336\", \"938\u00d79=8442\"),\n    @(\"797\u00d73=2391\", \"276\u00d74=1104\"),\n    @(\"362\u00d79=3258\", \"933\u00d77=6531\"),\n    @(\"431\u00d73=1293\", \"270\u00d73=810\"),\n    @(\"447\u00d75=2235\", \"533\u00d74=2132\"),\n    @(\"499\u00d76=2994\", \"181\u00d79=1629\"),\n    @(\"119\u00d73=357\", \"629\u00d76=3774\"),\n    @(\"970\u00d79=8730\", \"122\u00d76=732\"),\n    @(\"478\u00d74=1912\", \"718\u00d79=6462\"),\n    @(\"679\u00d78=5432\", \"917\u00d79=8253\"),\n    @(\"804\u00d72=1608\", \"505\u00d79=4545\"),\n    @(\"463\u00d73=1389\", \"431\u00d76=2586\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
